$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1912.25
$ws.Range("I40").Value = 1899.7142
$ws.Range("J40").Value = 2000
$ws.Range("K40").Value = 1899.7142
$ws.Range("L40").Value = 2000
$ws.Range("M40").Value = -1724.7142
$ws.Range("N40").Value = -2350
$ws.Range("H101").Value = 540
$ws.Range("I101").Value = 900
$ws.Range("K101").Value = 2700
$ws.Range("M101").Value = -1078
$ws.Range("H103").Value = 500
$ws.Range("I103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("M103").ClearContents()
$ws.Range("H111").Value = 1756.1765
$ws.Range("I111").Value = 1105.8462
$ws.Range("J111").Value = 3869.75
$ws.Range("K111").Value = 3317.5386
$ws.Range("L111").Value = 11609.25
$ws.Range("M111").Value = -250.5385999999999
$ws.Range("N111").Value = -17743.25
$ws.Range("H112").Value = 2232.0527
$ws.Range("J112").Value = 2312.125
$ws.Range("L112").Value = 6936.375
$ws.Range("N112").Value = -9152.375
$ws.Range("H113").Value = 4709.6665
$ws.Range("I113").Value = 4557.5557
$ws.Range("J113").Value = 5166
$ws.Range("K113").Value = 4557.5557
$ws.Range("L113").Value = 5166
$ws.Range("M113").Value = -1303.5557
$ws.Range("N113").Value = -11674
$ws.Range("H129").Value = 1867.1428
$ws.Range("I129").Value = 520
$ws.Range("K129").Value = 1560
$ws.Range("M129").Value = 3440
$ws.Range("H138").Value = 8087.933
$ws.Range("J138").Value = 7807.793
$ws.Range("L138").Value = 23423.379
$ws.Range("N138").Value = -33703.379

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15954.841
$ws.Range("I32").Value = 8587.861999999999
$ws.Range("J32").Value = 30197.666
$ws.Range("K32").Value = 8587.861999999999
$ws.Range("L32").Value = 30197.666
$ws.Range("M32").Value = -8300.861999999999
$ws.Range("N32").Value = -30771.666
$ws.Range("H63").Value = 6579.8
$ws.Range("J63").Value = 6316.5
$ws.Range("L63").Value = 6316.5
$ws.Range("N63").Value = -7688.5
$ws.Range("H66").Value = 6579.8
$ws.Range("J66").Value = 6316.5
$ws.Range("L66").Value = 31582.5
$ws.Range("N66").Value = -38446.5
$ws.Range("H74").Value = 2237.84
$ws.Range("I74").Value = 977.8095
$ws.Range("K74").Value = 977.8095
$ws.Range("M74").Value = -103.8095
$ws.Range("H77").Value = 2237.84
$ws.Range("I77").Value = 977.8095
$ws.Range("K77").Value = 4889.0475
$ws.Range("M77").Value = -521.0474999999997
$ws.Range("H95").Value = 54999.5
$ws.Range("J95").Value = 54999.5
$ws.Range("L95").Value = 54999.5
$ws.Range("N95").Value = -60491.5
$ws.Range("H122").Value = 773408.1
$ws.Range("I122").Value = 1432043.8
$ws.Range("K122").Value = 4296131.4
$ws.Range("M122").Value = -4293681.4
$ws.Range("H132").Value = 2304.3333
$ws.Range("I132").Value = 1492.8948
$ws.Range("K132").Value = 4478.6844
$ws.Range("M132").Value = -1948.6844

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 834.2
$ws.Range("I36").Value = 793.5
$ws.Range("J36").Value = 997
$ws.Range("K36").Value = 793.5
$ws.Range("L36").Value = 997
$ws.Range("M36").Value = -259.5
$ws.Range("N36").Value = -2065
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H134").Value = 2779.88
$ws.Range("I134").Value = 1425.9375
$ws.Range("K134").Value = 4277.8125
$ws.Range("M134").Value = -1742.8125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 119.94444
$ws.Range("I7").Value = 124.70588
$ws.Range("K7").Value = 124.70588
$ws.Range("M7").Value = -11.70587999999999
$ws.Range("H22").Value = 472.6
$ws.Range("I22").Value = 121.333336
$ws.Range("K22").Value = 121.333336
$ws.Range("M22").Value = 228.666664
$ws.Range("H31").Value = 5108.6523
$ws.Range("I31").Value = 3902.077
$ws.Range("J31").Value = 6677.2
$ws.Range("K31").Value = 3902.077
$ws.Range("L31").Value = 6677.2
$ws.Range("M31").Value = -3607.077
$ws.Range("N31").Value = -7267.2
$ws.Range("H34").Value = 5108.6523
$ws.Range("I34").Value = 3902.077
$ws.Range("J34").Value = 6677.2
$ws.Range("K34").Value = 3902.077
$ws.Range("L34").Value = 6677.2
$ws.Range("M34").Value = -3700.077
$ws.Range("N34").Value = -7081.2
$ws.Range("H58").Value = 3670.7058
$ws.Range("I58").Value = 1497.1666
$ws.Range("K58").Value = 1497.1666
$ws.Range("M58").Value = -1294.1666
$ws.Range("H99").Value = 17040
$ws.Range("I99").Value = 15498.25
$ws.Range("J99").Value = 17810.875
$ws.Range("K99").Value = 15498.25
$ws.Range("L99").Value = 17810.875
$ws.Range("M99").Value = -14000.25
$ws.Range("N99").Value = -20806.875
$ws.Range("H126").Value = 17040
$ws.Range("I126").Value = 15498.25
$ws.Range("J126").Value = 17810.875
$ws.Range("K126").Value = 46494.75
$ws.Range("L126").Value = 53432.625
$ws.Range("M126").Value = -44024.75
$ws.Range("N126").Value = -58372.625
$ws.Range("H134").Value = 2301.3157
$ws.Range("I134").Value = 1692.2142
$ws.Range("K134").Value = 5076.642599999999
$ws.Range("M134").Value = -2541.642599999999
$ws.Range("H136").Value = 3670.7058
$ws.Range("I136").Value = 1497.1666
$ws.Range("K136").Value = 4491.4998
$ws.Range("M136").Value = -1941.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 392.7
$ws.Range("I5").Value = 392.7
$ws.Range("K5").Value = 1178.1
$ws.Range("M5").Value = -1066.1
$ws.Range("H34").Value = 3616.6667
$ws.Range("J34").Value = 5000
$ws.Range("L34").Value = 15000
$ws.Range("N34").Value = -15168
$ws.Range("H39").Value = 1666.6666
$ws.Range("I39").Value = 1666.6666
$ws.Range("K39").Value = 4999.9998
$ws.Range("M39").Value = -4705.9998
$ws.Range("H46").Value = 834727.25
$ws.Range("I46").Value = 1021.25
$ws.Range("J46").Value = 1251580.2
$ws.Range("K46").Value = 3063.75
$ws.Range("L46").Value = 3754740.6
$ws.Range("M46").Value = -2972.75
$ws.Range("N46").Value = -3754922.6
$ws.Range("H51").Value = 597.6667
$ws.Range("I51").Value = 396.5
$ws.Range("K51").Value = 1189.5
$ws.Range("M51").Value = -729.5
$ws.Range("H55").Value = 68399.60000000001
$ws.Range("I55").Value = 200758.8
$ws.Range("J55").Value = 2220
$ws.Range("K55").Value = 602276.3999999999
$ws.Range("L55").Value = 6660
$ws.Range("M55").Value = -602099.3999999999
$ws.Range("N55").Value = -7014
$ws.Range("H113").Value = 3755.875
$ws.Range("I113").Value = 900
$ws.Range("J113").Value = 3946.2666
$ws.Range("K113").Value = 2700
$ws.Range("L113").Value = 11838.7998
$ws.Range("M113").Value = -530
$ws.Range("N113").Value = -16178.7998
$ws.Range("H135").Value = 392.7
$ws.Range("I135").Value = 392.7
$ws.Range("K135").Value = 3534.3
$ws.Range("M135").Value = -999.2999999999997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9665.223
$ws.Range("I80").Value = 9995
$ws.Range("J80").Value = 9624
$ws.Range("K80").Value = 9995
$ws.Range("L80").Value = 9624
$ws.Range("M80").Value = -8997
$ws.Range("N80").Value = -11620
$ws.Range("H83").Value = 9665.223
$ws.Range("I83").Value = 9995
$ws.Range("J83").Value = 9624
$ws.Range("K83").Value = 49975
$ws.Range("L83").Value = 48120
$ws.Range("M83").Value = -44983
$ws.Range("N83").Value = -58104
$ws.Range("H102").Value = 1699.9615
$ws.Range("I102").Value = 728.2353000000001
$ws.Range("K102").Value = 728.2353000000001
$ws.Range("M102").Value = 893.7646999999999
$ws.Range("H141").Value = 71935.5
$ws.Range("J141").Value = 71935.5
$ws.Range("L141").Value = 71935.5
$ws.Range("N141").Value = -82295.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2833.5
$ws.Range("I16").Value = 2897.6155
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 2897.6155
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -2727.6155
$ws.Range("N16").Value = -2340
$ws.Range("H46").Value = 3272.6365
$ws.Range("I46").Value = 2666.6667
$ws.Range("J46").Value = 3499.875
$ws.Range("K46").Value = 2666.6667
$ws.Range("L46").Value = 3499.875
$ws.Range("M46").Value = -2478.6667
$ws.Range("N46").Value = -3875.875
$ws.Range("H61").Value = 6842.3335
$ws.Range("I61").Value = 6010.3335
$ws.Range("K61").Value = 6010.3335
$ws.Range("M61").Value = -5808.3335
$ws.Range("H82").Value = 2789
$ws.Range("J82").Value = 2333
$ws.Range("L82").Value = 2333
$ws.Range("N82").Value = -3055
$ws.Range("H85").Value = 2789
$ws.Range("J85").Value = 2333
$ws.Range("L85").Value = 2333
$ws.Range("N85").Value = -4829
$ws.Range("H93").Value = 882.3333
$ws.Range("I93").Value = 558.3333
$ws.Range("K93").Value = 558.3333
$ws.Range("M93").Value = 689.6667
$ws.Range("H113").Value = 6842.3335
$ws.Range("I113").Value = 6010.3335
$ws.Range("K113").Value = 6010.3335
$ws.Range("M113").Value = -3840.3335
$ws.Range("H132").Value = 5630
$ws.Range("I132").Value = 3600
$ws.Range("K132").Value = 10800
$ws.Range("M132").Value = -8270

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1615.0526
$ws.Range("I132").Value = 707
$ws.Range("K132").Value = 2121
$ws.Range("M132").Value = 409

Write-Host "Applied all Seraph_Profits updates"